# Weekly Summary update for Wednesday July 17, 2013 entry.
#
# Strategy: this runtime auto-merges adjacent runs that share identical
# run-formatting as soon as any Range text edit touches them, so plain
# InsertBefore/InsertAfter calls can't reproduce the run-per-edit granularity
# seen in the target XML. However, Bookmarks.Add() onto a collapsed Range
# *does* force a permanent split of the surrounding run - and the split
# survives even after the bookmark itself is deleted. We use short-lived
# helper bookmarks purely to carve out the desired run boundaries, then
# remove them, leaving the real "_GoBack" bookmark seated in its proper
# (mid-word) position.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the paragraph ("Wednesday July 17, 2013 - ") and drop the
#    existing _GoBack bookmark (it will be re-created later at its new
#    position, further into the paragraph).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Last
$r = $p.Range
$paraStart = $r.Start

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Trim the trailing "- " down to a single trailing space, without
#    disturbing the run's leading <w:lastRenderedPageBreak/> -- we only
#    touch the very end of the run's text.
# ---------------------------------------------------------------------
$dashSpace = $d.Range($paraStart + 24, $paraStart + 26)
$dashSpace.Text = ""

# ---------------------------------------------------------------------
# 3. Insert the remaining new content right after "Wednesday July 17, 2013 ".
#    (All of it lands in one merged run for now; it gets carved into the
#    five target runs in step 4.)
# ---------------------------------------------------------------------
$insertPoint = $d.Range($paraStart + 24, $paraStart + 24)
$newTail = [char]0x2013 + " Had a meeting with Matt today. I received a yellow and need to work harder. Full description of the meeting is in the Tasks document."
$insertPoint.InsertBefore($newTail)

# ---------------------------------------------------------------------
# 4. Carve run boundaries using disposable bookmarks, matching the
#    target's run layout:
#      R1 "Wednesday July 17, 2013 "
#      R2 "\u2013"
#      R3 " "
#      R4 "Had a meeting ... Tasks docume"
#      [[_GoBack]]
#      R5 "nt."
# ---------------------------------------------------------------------
$b1 = $paraStart + 24                 # after "Wednesday July 17, 2013 "
$b2 = $b1 + 1                         # after the en dash
$b3 = $b2 + 1                         # after the following space
$b4 = $b3 + 130                       # after "...Tasks docume"

$d.Bookmarks.Add("_tmpSplit1", $d.Range($b1, $b1)) | Out-Null
$d.Bookmarks.Add("_tmpSplit2", $d.Range($b2, $b2)) | Out-Null
$d.Bookmarks.Add("_tmpSplit3", $d.Range($b3, $b3)) | Out-Null

# Re-seat _GoBack between "docume" and "nt." -- this both restores the
# bookmark and forces the R4 / R5 split.
$d.Bookmarks.Add("_GoBack", $d.Range($b4, $b4)) | Out-Null

$d.Bookmarks("_tmpSplit1").Delete()
$d.Bookmarks("_tmpSplit2").Delete()
$d.Bookmarks("_tmpSplit3").Delete()

# ---------------------------------------------------------------------
# 5. Append two new empty paragraphs after the entry.
# ---------------------------------------------------------------------
$end = $d.Paragraphs.Last.Range
$end.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertParagraphAfter()
